$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A102").Value = "Record"
$ws.Range("B102").Value = "RJ Record"
$ws.Range("C102").Value = "Economia"
$ws.Range("D102").Value = "2025-04-09T18:24"
$ws.Range("E102").Value = "Positivo"
$ws.Range("F102").Value = "Imposto de Renda. Mutirão no Centro de Campos tira dúvidas sobre preenchimento da declaração. Entrevista com senhora que foi atendida e com o subsecretário de Desenvolvimento Econômico, Felipe Knust. Mutirão foi hoje e muita gente aproveitou a oportunidade. Evento é parceria entre o Conselho Regional de Contabilidade do RJ e a Prefeitura de Campos. Equipe de contadores voluntários atendendo. Entrevista com delegada do CRC, Fabiana Viana. *matéria*"

$ws.Range("A103").Value = "Record"
$ws.Range("B103").Value = "RJ Record"
$ws.Range("C103").Value = "Defesa Civil"
$ws.Range("D103").Value = "2025-04-09T19:11"
$ws.Range("E103").Value = "Positivo"
$ws.Range("F103").Value = "Disposofobia. Transtorno de acumulação compulsiva que afeta cerca de 4% da população mundial. Entrevista com psicanalista, Arthur Costa.  Esta semana, um caso aconteceu em Cabo Frio. Recentemente, também teve um caso em Campos. Defesa Civil teve que fazer poda nas árvores na ocasião para acessar a casa do Parque Bela Vista. *matéria* também foi exibida no Balanço Geral."
